$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily record for 2026/01/25 07:00 (rank 19, score 139) needs to be
# inserted into the log at row 723, pushing the existing rows 723-764 down
# to 724-765 (the sheet's dimension grows from D764 to D765).
$ws.Rows.Item(723).Insert()

# Column A stores dates as plain text (e.g. "2026/01/25"), not Excel date
# serials. Force the new cell to a text format before assigning the value
# so it isn't auto-converted to a date, then restore the default "Normal"
# style so the cell matches the rest of the column (no explicit style).
$ws.Range("A723").NumberFormat = "@"
$ws.Range("A723").Value = "2026/01/25"
$ws.Range("A723").Style = "Normal"

$ws.Range("B723").Value = "日"
$ws.Range("C723").Value = 19
$ws.Range("D723").Value = 139
